$wb = $excel.ActiveWorkbook

# --- Insert new sheet "2022-Q4" right after the "总计" sheet (position 2).
#     We duplicate the existing "2022-Q3" sheet (Worksheet.Copy keeps all
#     cell formatting/styles) and then overwrite its cell values - this
#     avoids cross-sheet Copy/PasteSpecial which does not transfer
#     formatting in this runtime. ---
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q4"

# The source sheet only has one data row; duplicate it (same-sheet copy
# preserves formatting) so we have two styled data rows to work with.
$new.Range("A2:H2").Copy($new.Range("A3:H3"))

# --- Header row (values identical to other sheets, just re-assert) ---
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Force text storage for the fund columns (matches source data which is
# stored as text, not numeric, in the original workbook), then drop back
# to the default/unstyled cell format (style index 0) so the cells carry
# no explicit style, exactly as in the other quarter sheets.
$new.Range("B2:G3").NumberFormat = "@"

# --- Data rows ---
$new.Range("A2").Value = 0
$new.Range("B2").Value = "159743"
$new.Range("C2").Value = "博时中证湖北新旧动能转换ETF"
$new.Range("D2").Value = "3.43"
$new.Range("E2").Value = "99.18"
$new.Range("F2").Value = "1.66"
$new.Range("G2").Value = "0.0569"
$new.Range("H2").Value = 9

$new.Range("A3").Value = 1
$new.Range("B3").Value = "519677"
$new.Range("C3").Value = "银河定投宝中证腾讯济安价值100A股指数"
$new.Range("D3").Value = "3.05"
$new.Range("E3").Value = "91.25"
$new.Range("F3").Value = "1.21"
$new.Range("G3").Value = "0.0369"
$new.Range("H3").Value = 7

$new.Range("B2:G3").Style = "Normal"

# --- Update the "总计" summary sheet: insert a new row for 2022-Q4 above
#     the existing 2022-Q3 row, shifting the rest down. ---
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.09

# Copy style (s="2") to the new A2 cell from A3 (same style, same sheet).
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0

# Renumber the index column (A) for the remaining rows to stay 0-based
# sequential, matching the sheet's existing convention.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

Write-Output "done"
